# Updates the cryptos price/volume columns (D, E) for rows 2-51 and
# swaps the Toncoin/RenderToken rows (29-30), matching the latest
# GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '62.054.63'
$ws.Range("E2").Value = '  -0.96%  '

# Row 3
$ws.Range("D3").Value = '3.414.26'
$ws.Range("E3").Value = '  -0.85%  '

# Row 4
$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.23%  '

# Row 5
$ws.Range("D5").Value = '''410.63'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.69%  '

# Row 6
$ws.Range("D6").Value = '''129.19'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.49%  '

# Row 7
$ws.Range("D7").Value = '''0.641'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +6.34%  '

# Row 8
$ws.Range("D8").Value = '''0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.02%  '

# Row 9
$ws.Range("E9").Value = '  +5.16%  '

# Row 10
$ws.Range("E10").Value = '  +1.10%  '

# Row 11
$ws.Range("D11").Value = '''43.72'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.95%  '

# Row 12
$ws.Range("D12").Value = '''0.0000228'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +40.36%  '

# Row 13
$ws.Range("E13").Value = '  +9.34%  '

# Row 14
$ws.Range("E14").Value = '  -0.27%  '

# Row 15
$ws.Range("D15").Value = '''21.39'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.86%  '

# Row 16
$ws.Range("D16").Value = '3.954.08'
$ws.Range("E16").Value = '  -0.86%  '

# Row 17
$ws.Range("D17").Value = '3.444.81'
$ws.Range("E17").Value = '  +0.60%  '

# Row 18
$ws.Range("E18").Value = '  +8.58%  '

# Row 19
$ws.Range("E19").Value = '  +6.02%  '

# Row 20
$ws.Range("D20").Value = '61.974.49'
$ws.Range("E20").Value = '  -1.16%  '

# Row 21
$ws.Range("D21").Value = '''477.92'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +51.12%  '

# Row 22
$ws.Range("D22").Value = '''93.11'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +9.99%  '

# Row 23
$ws.Range("D23").Value = '''3.22'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.37%  '

# Row 24
$ws.Range("D24").Value = '''13.24'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.41%  '

# Row 25
$ws.Range("E25").Value = '  +4.11%  '

# Row 26
$ws.Range("D26").Value = '''33.54'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +11.22%  '

# Row 27
$ws.Range("D27").Value = '''9.35'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +14.10%  '

# Row 28
$ws.Range("D28").Value = '''4.76'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.09%  '

# Row 29
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = '''7.64'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.98%  '

# Row 30
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = '''2.71'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.22%  '

# Row 31
$ws.Range("D31").Value = '''12.11'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.92%  '

# Row 32
$ws.Range("E32").Value = '  -2.26%  '

# Row 33
$ws.Range("D33").Value = '''0.115'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.15%  '

# Row 34
$ws.Range("D34").Value = '''42.45'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.36%  '

# Row 35
$ws.Range("E35").Value = '  +0.00%  '

# Row 36
$ws.Range("D36").Value = '''0.0505'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.56%  '

# Row 37
$ws.Range("D37").Value = '''53.92'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.90%  '

# Row 38
$ws.Range("D38").Value = '''0.998'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.06%  '

# Row 39
$ws.Range("E39").Value = '  +7.90%  '

# Row 40
$ws.Range("E40").Value = '  +1.91%  '

# Row 41
$ws.Range("E41").Value = '  -0.76%  '

# Row 42
$ws.Range("D42").Value = '''0.319'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.16%  '

# Row 43
$ws.Range("D43").Value = '''4.41'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +11.71%  '

# Row 44
$ws.Range("D44").Value = '''144.50'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.62%  '

# Row 45
$ws.Range("E45").Value = '  +16.33%  '

# Row 46
$ws.Range("E46").Value = '  +1.37%  '

# Row 47
$ws.Range("D47").Value = '''16.71'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.07%  '

# Row 48
$ws.Range("D48").Value = '''0.150'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +19.98%  '

# Row 49
$ws.Range("D49").Value = '''22.62'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.79%  '

# Row 50
$ws.Range("D50").Value = '''2.16'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.40%  '

# Row 51
$ws.Range("D51").Value = '3.753.65'
$ws.Range("E51").Value = '  -0.99%  '
